# Scheduled price-refresh run: pushes freshly-pulled market-board averages
# (currentAveragePrice / NQ / HQ) and the recomputed Leve price/profit
# columns (H,I,J,K,L,M,N) into each crafting-class sheet. Values below were
# produced by the external pricing job; this script only writes them.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2570.3
$ws.Range("I8").Value = 50.75
$ws.Range("K8").Value = 152.25
$ws.Range("M8").Value = -13.25

$ws.Range("H43").Value = 2750
$ws.Range("I43").Value = 2750
$ws.Range("K43").Value = 2750
$ws.Range("M43").Value = -2681

$ws.Range("H64").Value = 5522.5713
$ws.Range("I64").Value = 4531.6
$ws.Range("J64").Value = 8000
$ws.Range("K64").Value = 4531.6
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = -4283.6
$ws.Range("N64").Value = -8496

$ws.Range("H67").Value = 5522.5713
$ws.Range("I67").Value = 4531.6
$ws.Range("J67").Value = 8000
$ws.Range("K67").Value = 4531.6
$ws.Range("L67").Value = 8000
$ws.Range("M67").Value = -3673.6
$ws.Range("N67").Value = -9716

$ws.Range("H100").Value = 3672.138
$ws.Range("I100").Value = 1593.25
$ws.Range("K100").Value = 1593.25
$ws.Range("M100").Value = -1052.25

$ws.Range("H136").Value = 88070
$ws.Range("J136").Value = 88070
$ws.Range("L136").Value = 88070
$ws.Range("N136").Value = -98270

$ws.Range("H137").Value = 20006242
$ws.Range("J137").Value = 14482.375
$ws.Range("L137").Value = 43447.125
$ws.Range("N137").Value = -48547.125

$ws.Range("H139").Value = 74780
$ws.Range("J139").Value = 74780
$ws.Range("L139").Value = 74780
$ws.Range("N139").Value = -85060

$ws.Range("H140").Value = 63380
$ws.Range("J140").Value = 76725
$ws.Range("L140").Value = 76725
$ws.Range("N140").Value = -87085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1012969.8

$ws.Range("H45").Value = 1623.6
$ws.Range("I45").Value = 1486.3684
$ws.Range("K45").Value = 1486.3684
$ws.Range("M45").Value = -1109.3684

$ws.Range("H61").Value = 2506011.2
$ws.Range("I61").Value = 7075.346
$ws.Range("K61").Value = 7075.346
$ws.Range("M61").Value = -6863.346

$ws.Range("H132").Value = 5840.5356
$ws.Range("I132").Value = 4186.75
$ws.Range("J132").Value = 7080.875
$ws.Range("K132").Value = 12560.25
$ws.Range("L132").Value = 21242.625
$ws.Range("M132").Value = -10030.25
$ws.Range("N132").Value = -26302.625

$ws.Range("H136").Value = 2506011.2
$ws.Range("I136").Value = 7075.346
$ws.Range("K136").Value = 21226.038
$ws.Range("M136").Value = -18676.038

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6963.9473
$ws.Range("I99").Value = 9199.076999999999
$ws.Range("K99").Value = 9199.076999999999
$ws.Range("M99").Value = -7701.076999999999

$ws.Range("H107").Value = 1619.9333
$ws.Range("I107").Value = 1394.3889
$ws.Range("J107").Value = 1958.25
$ws.Range("K107").Value = 1394.3889
$ws.Range("L107").Value = 1958.25
$ws.Range("M107").Value = 525.6111000000001
$ws.Range("N107").Value = -5798.25

$ws.Range("H134").Value = 3973559.2
$ws.Range("I134").Value = 6855.533
$ws.Range("K134").Value = 20566.599
$ws.Range("M134").Value = -18031.599

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10741440
$ws.Range("I58").Value = 27780412
$ws.Range("J58").Value = 3925850.2
$ws.Range("K58").Value = 27780412
$ws.Range("L58").Value = 3925850.2
$ws.Range("M58").Value = -27780209
$ws.Range("N58").Value = -3926256.2

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H107").Value = 783.1875
$ws.Range("I107").Value = 464.27274
$ws.Range("J107").Value = 1484.8
$ws.Range("K107").Value = 464.27274
$ws.Range("L107").Value = 1484.8
$ws.Range("M107").Value = 1455.72726
$ws.Range("N107").Value = -5324.8

$ws.Range("H122").Value = 13659.137
$ws.Range("I122").Value = 2215.8333
$ws.Range("J122").Value = 27391.1
$ws.Range("K122").Value = 6647.499899999999
$ws.Range("L122").Value = 82173.29999999999
$ws.Range("M122").Value = -4197.499899999999
$ws.Range("N122").Value = -87073.29999999999

$ws.Range("H136").Value = 10741440
$ws.Range("I136").Value = 27780412
$ws.Range("J136").Value = 3925850.2
$ws.Range("K136").Value = 83341236
$ws.Range("L136").Value = 11777550.6
$ws.Range("M136").Value = -83338686
$ws.Range("N136").Value = -11782650.6

$ws.Range("H140").Value = 78400
$ws.Range("J140").Value = 78400
$ws.Range("L140").Value = 78400
$ws.Range("N140").Value = -88760

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4035970
$ws.Range("I122").Value = 8066940.5
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 72602464.5
$ws.Range("L122").Value = 44995.5
$ws.Range("M122").Value = -72600014.5
$ws.Range("N122").Value = -49895.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 976.4
$ws.Range("I107").Value = 1070.5
$ws.Range("K107").Value = 1070.5
$ws.Range("M107").Value = 849.5

$ws.Range("H113").Value = 1499.9231
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 47514
$ws.Range("I42").Value = 35000
$ws.Range("K42").Value = 35000
$ws.Range("M42").Value = -34437

$ws.Range("H49").Value = 47514
$ws.Range("I49").Value = 35000
$ws.Range("K49").Value = 35000
$ws.Range("M49").Value = -34853

$ws.Range("H132").Value = 1114498.9
$ws.Range("I132").Value = 2567269.5
$ws.Range("J132").Value = 3556.4707
$ws.Range("K132").Value = 7701808.5
$ws.Range("L132").Value = 10669.4121
$ws.Range("M132").Value = -7699278.5
$ws.Range("N132").Value = -15729.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 37142.715
$ws.Range("I33").Value = 20000
$ws.Range("J33").Value = 39999.832
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 39999.832
$ws.Range("M33").Value = -19750
$ws.Range("N33").Value = -40499.832

$ws.Range("H36").Value = 37142.715
$ws.Range("I36").Value = 20000
$ws.Range("J36").Value = 39999.832
$ws.Range("K36").Value = 20000
$ws.Range("L36").Value = 39999.832
$ws.Range("M36").Value = -19750
$ws.Range("N36").Value = -40499.832

$ws.Range("H45").Value = 45000
$ws.Range("J45").Value = 45000
$ws.Range("L45").Value = 45000
$ws.Range("N45").Value = -45982

$ws.Range("H113").Value = 1769.1212
$ws.Range("I113").Value = 1493.6111
$ws.Range("K113").Value = 4480.8333
$ws.Range("M113").Value = -2310.8333

$ws.Range("H136").Value = 5051417.5
$ws.Range("I136").Value = 2420933
$ws.Range("J136").Value = 20834322
$ws.Range("K136").Value = 7262799
$ws.Range("L136").Value = 62502966
$ws.Range("M136").Value = -7260249
$ws.Range("N136").Value = -62508066
